# Update the "F" column (想去人数 / want-to-go count) values on the
# "展览" and "全部类型" worksheets to match the regenerated page output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# row -> new value for sheet "展览"
$sheet1Updates = @{
    3  = 221
    4  = 0
    5  = 206
    6  = 146
    7  = 0
    9  = 91
    11 = 0
    12 = 1128
    13 = 99
    14 = 0
    15 = 0
    16 = 80
    17 = 0
    18 = 109
    19 = 0
    20 = 0
    22 = 38
    23 = 0
    24 = 533
    26 = 3945
    27 = 394
    30 = 568
    31 = 526
    32 = 137
    34 = 294
    36 = 0
    37 = 0
    38 = 942
    40 = 0
    41 = 56
    42 = 483
    44 = 0
}

# row -> new value for sheet "全部类型"
$sheet4Updates = @{
    2  = 29
    3  = 0
    5  = 206
    6  = 0
    7  = 110
    8  = 0
    12 = 215
    14 = 99
    17 = 80
    18 = 140
    19 = 109
    20 = 3853
    21 = 6190
    22 = 0
    23 = 0
    24 = 0
    25 = 533
    26 = 0
    27 = 0
    29 = 0
    30 = 2541
    32 = 526
    35 = 294
    36 = 363
    37 = 0
    38 = 1554
    39 = 942
    40 = 42
    41 = 52
    42 = 56
    43 = 483
    45 = 73
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
